$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.093.15'
$ws.Range("D3").Value = '1.790.67'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''226.98'
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("D6").Value = '''0.547'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D8").Value = '''32.37'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +3.87%  '
$ws.Range("D10").Value = '''0.0690'
$ws.Range("E10").Value = '  -2.07%  '
$ws.Range("D11").Value = '''0.0939'
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Value = '2.048.25'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '''11.45'
$ws.Range("E13").Value = '  +5.98%  '
$ws.Range("D14").Value = '1.791.23'
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").Value = '''0.624'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '34.076.41'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").Value = '''68.10'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '''243.99'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").Value = '0.0₃0783'
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").Value = '''10.95'
$ws.Range("E21").Value = '  +1.86%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").Value = '''2.05'
$ws.Range("E24").Value = '  -2.67%  '
$ws.Range("D25").Value = '''162.07'
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("E26").Value = '  +2.57%  '
$ws.Range("D27").Value = '''16.30'
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +1.50%  '
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("D31").Value = '''0.0521'
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("D32").Value = '''3.67'
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D34").Value = '''1.85'
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("D35").Value = '1.414.80'
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("E37").Value = '  +2.79%  '
$ws.Range("E38").Value = '  +7.50%  '
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("E40").Value = '  +3.28%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''0.923'
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '''2.35'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '''13.54'
$ws.Range("E44").Value = '  +9.39%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '''6.08'
$ws.Range("E45").Value = '  +3.63%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").Value = '''0.0507'
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0138'
$ws.Range("E47").Value = '  -4.11%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = '''107.38'
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").Value = '1.949.34'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("E51").Value = '  +0.06%  '
